$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1026
$ws1.Range("F4").Value = 166
$ws1.Range("F5").Value = 2773
$ws1.Range("F7").Value = 218
$ws1.Range("F10").Value = 60
$ws1.Range("F11").Value = 66
$ws1.Range("F12").Value = 2587
$ws1.Range("F13").Value = 753

# Sheet "全部类型" (sheet4): update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1026
$ws4.Range("F5").Value = 166
$ws4.Range("F6").Value = 2773
$ws4.Range("F8").Value = 218
$ws4.Range("F12").Value = 60
$ws4.Range("F13").Value = 66
$ws4.Range("F14").Value = 2587
$ws4.Range("F15").Value = 753
